$p = $ppt.ActivePresentation

# Slide 4 ("For Live Session: Question 3") - Content Placeholder 2
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# --- Edit 1 -------------------------------------------------------------
# "...regressors EXCEPT for mps (since..." ->
# "...regressors EXCEPT for mpg (since..."
# (the word that used to be flagged as a misspelling, "mps", becomes "mpg";
#  "for " moves out of the first run and the leading space moves out of the
#  trailing run, so replace the whole "for mps " span in one shot)
$full = $tr.Text
$needle1 = "for mps "
$idx1 = $full.IndexOf($needle1)
if ($idx1 -ge 0) {
    $rng1 = $tr.Characters($idx1 + 1, $needle1.Length)
    $rng1.Text = "for mpg "
}

# --- Edit 2 -------------------------------------------------------------
# "Assess the relationship between the mpg and the slope." ->
# "Assess the relationship between the mpg and the horsepower."
$full = $tr.Text
$needle2 = "Assess the relationship between the mpg and the slope.  Make sure and include estimates of your uncertainty ("
$idx2 = $full.IndexOf($needle2)
if ($idx2 -ge 0) {
    $rng2 = $tr.Characters($idx2 + 1, $needle2.Length)
    $rng2.Text = "Assess the relationship between the mpg and the horsepower.  Make sure and include estimates of your uncertainty ("
}
